# Updates the cryptos price list: refreshed Price/Volume(1h) figures for
# most rows, plus the Toncoin/Monero rows (29-30) swapping places with new
# coin/link/price/volume data.
# Commit: "Updated cryptos list on Fri Dec 15 18:51:55 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.110.04"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "2.242.08"
$ws.Range("E3").Value = "  -2.19%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'247.51"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").Value = "'0.624"
$ws.Range("E6").Value = "  -2.54%  "
$ws.Range("D7").Value = "'77.12"
$ws.Range("E7").Value = "  +4.28%  "
$ws.Range("D9").Value = "'0.631"
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("D10").Value = "'41.17"
$ws.Range("E10").Value = "  +4.90%  "
$ws.Range("D11").Value = "'0.0957"
$ws.Range("E11").Value = "  -2.68%  "
$ws.Range("D12").Value = "'7.19"
$ws.Range("E12").Value = "  -3.39%  "
$ws.Range("E13").Value = "  -2.92%  "
$ws.Range("D14").Value = "2.575.99"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").Value = "'14.87"
$ws.Range("E15").Value = "  -2.56%  "
$ws.Range("D16").Value = "'0.859"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "2.243.70"
$ws.Range("E17").Value = "  -2.78%  "
$ws.Range("D18").Value = "41.932.57"
$ws.Range("E18").Value = "  -2.43%  "
$ws.Range("D19").Value = "0.0₃0983"
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").Value = "'6.12"
$ws.Range("E20").Value = "  -3.33%  "
$ws.Range("D21").Value = "'71.82"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("E22").Value = "  +5.06%  "
$ws.Range("D23").Value = "'231.68"
$ws.Range("E23").Value = "  -2.56%  "
$ws.Range("D24").Value = "'11.44"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -5.63%  "
$ws.Range("E27").Value = "  -5.07%  "
$ws.Range("D28").Value = "'7.27"
$ws.Range("E28").Value = "  +12.49%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'168.81"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.10"
$ws.Range("E30").Value = "  -4.23%  "
$ws.Range("D31").Value = "'20.57"
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("D32").Value = "'33.04"
$ws.Range("E32").Value = "  +6.09%  "
$ws.Range("E33").Value = "  +1.54%  "
$ws.Range("E34").Value = "  -5.36%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "'4.50"
$ws.Range("E36").Value = "  -2.40%  "
$ws.Range("D37").Value = "'4.90"
$ws.Range("E37").Value = "  +2.70%  "
$ws.Range("D38").Value = "'0.0302"
$ws.Range("E38").Value = "  -2.58%  "
$ws.Range("D39").Value = "'14.23"
$ws.Range("E39").Value = "  -1.43%  "
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("E41").Value = "  -6.86%  "
$ws.Range("D42").Value = "'113.29"
$ws.Range("E42").Value = "  +14.06%  "
$ws.Range("E43").Value = "  -6.64%  "
$ws.Range("D44").Value = "'61.43"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("D45").Value = "'8.73"
$ws.Range("E45").Value = "  -4.03%  "
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("D47").Value = "'0.996"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("E48").Value = "  -3.06%  "
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").Value = "'4.27"
$ws.Range("E50").Value = "  -13.21%  "
$ws.Range("D51").Value = "'2.27"
$ws.Range("E51").Value = "  -2.13%  "
